$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""69.371.43"""
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Formula = "=""3.429.27"""
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").Formula = "=""0.999"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Formula = "=""580.46"""
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").Formula = "=""177.08"""
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Formula = "=""3.421.01"""
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Formula = "=""1.00"""
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Formula = "=""0.585"""
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Formula = "=""48.66"""
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Formula = "=""0.0000281"""
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Formula = "=""698.17"""
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Formula = "=""3.976.39"""
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Formula = "=""8.64"""
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Formula = "=""69.422.13"""
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Formula = "=""3.425.92"""
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Formula = "=""17.75"""
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Formula = "=""11.45"""
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Formula = "=""5.40"""
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Formula = "=""16.98"""
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").Formula = "=""9.63"""
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Formula = "=""33.65"""
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Formula = "=""8.78"""
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").Formula = "=""7.00"""
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Formula = "=""3.80"""
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("D33").Formula = "=""569.58"""
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("D36").Formula = "=""58.15"""
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Formula = "=""1.00"""
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Formula = "=""3.612.58"""
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -2.99%  "
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("D40").Formula = "=""34.94"""
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Formula = "=""0.0₃0735"""
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("D42").Formula = "=""3.29"""
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").Formula = "=""2.68"""
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Formula = "=""3.35"""
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Formula = "=""0.334"""
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Formula = "=""0.0420"""
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Formula = "=""1.48"""
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +4.83%  "
$ws.Range("D48").Formula = "=""2.67"""
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Formula = "=""0.128"""
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").Formula = "=""0.999"""
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Formula = "=""131.21"""
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  -1.20%  "
$excel.CutCopyMode = 0
